$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$range = $ws.Range("B2:E13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $excel.WorksheetFunction.Round([double]$val, 0)
    }
}
